$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell content updates on Sheet1 ---

# D11 used to show "TOTAL TEST CASES - 28"'s neighboring text; it should now
# read "Verify the functionality on home page."
$ws.Range("D11").Value2 = "Verify the functionality on home page."

# E16 test-case count bumped from 5 to 6
$ws.Range("E16").Value2 = 6

# E18 total test-case label bumped from 28 to 29
$ws.Range("E18").Value2 = "TOTAL TEST CASES - 29"

# Move the active selection from D11 to H10
$ws.Range("H10").Select()
